# Append a new transaction row (row 68) to the Transactions sheet.
# Every column in this sheet (including the numeric-looking ones) is
# stored as text, so force a text number format before writing the
# values - otherwise Excel would auto-coerce the date into a date
# serial number and the "Number of Units" / "Price per Share" columns
# into real numbers instead of text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A68:E68").NumberFormat = "@"

$ws.Range("A68").Value = "2025-09-22"
$ws.Range("B68").Value = "NESTLE"
$ws.Range("C68").Value = "Buy"
$ws.Range("D68").Value = "1"
$ws.Range("E68").Value = "5000"
